# Add Q4-2022 holdings data: insert a new "2022-Q4" sheet (per-fund holding
# detail, same layout as the other quarter sheets) right before the
# existing "2022-Q3" sheet, and prepend a matching summary row to the
# "总计" totals sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Structural changes first (sheet add + row insert shift cell
#    positions, so do them before grabbing any cell references).
# ---------------------------------------------------------------------
$q3ForAdd = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($q3ForAdd)
$q4.Name = "2022-Q4"

$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
# Insert() copies the formatting of the row above (the bold header row)
# into every column of the freshly inserted row; only column A should
# stay bold in the data rows, so strip the rest back to the default
# (un-styled) look before writing real content into it.
$total.Range("B2:D2").ClearFormats()

# ---------------------------------------------------------------------
# 2. Re-fetch a clean reference to "2022-Q3" (its position moved down
#    by one when the new sheet was inserted before it) and use it as a
#    formatting template for the new sheet - copy/paste-format keeps us
#    on the exact same shared style entries the rest of the workbook
#    already uses instead of re-deriving a near-duplicate style.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")

$q3.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$q3.Range("A2").Copy()
$q4.Range("A2:A10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row text
$q4.Cells.Item(1,2).Value = "基金代码"
$q4.Cells.Item(1,3).Value = "基金名称"
$q4.Cells.Item(1,4).Value = "基金规模"
$q4.Cells.Item(1,5).Value = "股票总仓位"
$q4.Cells.Item(1,6).Value = "仓位占比"
$q4.Cells.Item(1,7).Value = "持有市值(亿元)"
$q4.Cells.Item(1,8).Value = "仓位排名"

function Set-TextCell($ws, $row, $col, $text) {
    # Force text storage even for numeric-looking strings (fund codes,
    # percentages, etc. are stored as text in every sheet of this
    # workbook) by setting an explicit text format right before the
    # assignment.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

$q4Rows = @(
    @("006049", "恒越研究精选混合A/B",  "4.23", "88.62", "2.89", "0.1222", 9),
    @("007192", "恒越研究精选混合C",     "3.19", "88.62", "2.89", "0.0922", 9),
    @("014126", "华夏中证1000指数增强C", "8.66", "92.31", "0.87", "0.0753", 3),
    @("015963", "汇安品质优选混合A",     "2.48", "86.98", "2.91", "0.0722", 9),
    @("015635", "汇安价值先锋混合A",     "1.41", "84.53", "2.99", "0.0422", 9),
    @("015964", "汇安品质优选混合C",     "0.70", "86.98", "2.91", "0.0204", 9),
    @("350007", "天治趋势精选混合",      "0.39", "93.83", "3.68", "0.0144", 6),
    @("014125", "华夏中证1000指数增强A", "1.04", "92.31", "0.87", "0.0090", 3),
    @("015636", "汇安价值先锋混合C",     "0.09", "84.53", "2.99", "0.0027", 9)
)

$r = 2
foreach ($row in $q4Rows) {
    $q4.Cells.Item($r, 1).Value = $r - 2

    Set-TextCell $q4 $r 2 $row[0]
    Set-TextCell $q4 $r 3 $row[1]
    Set-TextCell $q4 $r 4 $row[2]
    Set-TextCell $q4 $r 5 $row[3]
    Set-TextCell $q4 $r 6 $row[4]
    Set-TextCell $q4 $r 7 $row[5]
    $q4.Cells.Item($r, 8).Value = $row[6]

    $r = $r + 1
}

# ---------------------------------------------------------------------
# 3. Fill in the new "总计" row. Row 2 was just inserted (blank, no
#    format) above the previous "2022-Q3" row, which is now row 3.
# ---------------------------------------------------------------------
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 9
$total.Cells.Item(2, 4).Value = 0.45

# The rows that got pushed down keep their old index-column numbering
# (0,1,2,3 again) - renumber them to 1,2,3,4 so the sequence stays
# contiguous under the new row.
for ($row = 3; $row -le 6; $row++) {
    $total.Cells.Item($row, 1).Value = $row - 2
}
